$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '39.630.38'
$ws.Range('E2').Value = '  +0.65%  '
$ws.Range('D3').Value = '2.169.33'
$ws.Range('E3').Value = '  +0.56%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range('D5').Value = '226.66'
$ws.Range("D5").Style = "Normal"
$ws.Range('E5').Value = '  -1.16%  '
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range('D7').Value = '63.11'
$ws.Range("D7").Style = "Normal"
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range('D10').Value = '0.0852'
$ws.Range("D10").Style = "Normal"
$ws.Range('E10').Value = '  -0.77%  '
$ws.Range('E12').Value = '  -1.87%  '
$ws.Range('D13').Value = '2.489.76'
$ws.Range('E13').Value = '  +0.36%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range('D14').Value = '21.76'
$ws.Range("D14").Style = "Normal"
$ws.Range('E14').Value = '  -2.34%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range('D15').Value = '0.809'
$ws.Range("D15").Style = "Normal"
$ws.Range('E15').Value = '  -0.94%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range('D16').Value = '5.49'
$ws.Range("D16").Style = "Normal"
$ws.Range('E16').Value = '  -1.23%  '
$ws.Range('D17').Value = '2.167.28'
$ws.Range('E17').Value = '  +0.52%  '
$ws.Range('D18').Value = '39.588.20'
$ws.Range('E18').Value = '  +0.10%  '
$ws.Range('D19').Value = '0.0₃0917'
$ws.Range('E19').Value = '  +7.55%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range('D20').Value = '71.74'
$ws.Range("D20").Style = "Normal"
$ws.Range('E20').Value = '  -0.96%  '
$ws.Range('E21').Value = '  -2.36%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range('D22').Value = '227.78'
$ws.Range("D22").Style = "Normal"
$ws.Range('E22').Value = '  -0.57%  '
$ws.Range('E24').Value = '  -3.79%  '
$ws.Range('E25').Value = '  -1.33%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range('D26').Value = '170.87'
$ws.Range("D26").Style = "Normal"
$ws.Range('E26').Value = '  -1.28%  '
$ws.Range('E27').Value = '  -2.13%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range('D28').Value = '0.139'
$ws.Range("D28").Style = "Normal"
$ws.Range('E28').Value = '  +0.65%  '
$ws.Range('E29').Value = '  +1.25%  '
$ws.Range('E30').Value = '  +0.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range('D31').Value = '2.67'
$ws.Range("D31").Style = "Normal"
$ws.Range('E31').Value = '  +3.76%  '
$ws.Range('E32').Value = '  +0.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range('D33').Value = '4.52'
$ws.Range("D33").Style = "Normal"
$ws.Range('E33').Value = '  -2.97%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range('D34').Value = '4.71'
$ws.Range("D34").Style = "Normal"
$ws.Range('E34').Value = '  -2.58%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range('D35').Value = '6.97'
$ws.Range("D35").Style = "Normal"
$ws.Range('E35').Value = '  -2.96%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range('D37').Value = '3.81'
$ws.Range("D37").Style = "Normal"
$ws.Range('E37').Value = '  +6.32%  '
$ws.Range('E38').Value = '  -0.95%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range('D39').Value = '0.999'
$ws.Range("D39").Style = "Normal"
$ws.Range('E39').Value = '  -0.26%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range('D40').Value = '4.97'
$ws.Range("D40").Style = "Normal"
$ws.Range('E40').Value = '  +19.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range('D41').Value = '102.48'
$ws.Range("D41").Style = "Normal"
$ws.Range('E41').Value = '  -0.81%  '
$ws.Range('E42').Value = '  -1.34%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range('D43').Value = '17.68'
$ws.Range("D43").Style = "Normal"
$ws.Range('D44').Value = '1.514.48'
$ws.Range('E44').Value = '  -1.70%  '
$ws.Range('E45').Value = '  +1.08%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range('D46').Value = '7.88'
$ws.Range("D46").Style = "Normal"
$ws.Range('E46').Value = '  +0.96%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").NumberFormat = "@"
$ws.Range('D47').Value = '0.0922'
$ws.Range("D47").Style = "Normal"
$ws.Range('E47').Value = '  -0.45%  '
$ws.Range('B48').Value = 'HuobiToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D48").NumberFormat = "@"
$ws.Range('D48').Value = '2.80'
$ws.Range("D48").Style = "Normal"
$ws.Range('E48').Value = '  +0.01%  '
$ws.Range('E49').Value = '  -1.66%  '
$ws.Range('E50').Value = '  +32.52%  '
$ws.Range('D51').Value = '2.369.42'
$ws.Range('E51').Value = '  +0.20%  '
